$wb = $excel.ActiveWorkbook

# --- "info" sheet: update row 1 values (kept as text, not numbers) ---
$info = $wb.Worksheets.Item("info")
$info.Range("A1:C1").NumberFormat = "@"
$info.Range("A1").Value = "111"
$info.Range("B1").Value = "111"
$info.Range("C1").Value = "1"
# D1 ("1") is unchanged.

# --- "items" sheet: update row 1, then drop rows 2 and 3 ---
$items = $wb.Worksheets.Item("items")
$items.Range("A1").Value = "녹말요지"
$items.Range("C1").Value = 1100
$items.Range("E1").Value = 1100
# B1 ("개") and D1 (1) are unchanged.

# Remove the now-unwanted rows 2 ("밥공기(1)") and 3 ("접시100").
# Deleting row 2 twice shifts row 3 up into row 2, then removes it too.
$items.Rows(2).Delete()
$items.Rows(2).Delete()
